$d = $word.ActiveDocument

$replacements = @(
    @("91÷7=13, 0", "52÷4=13, 0"),
    @("52÷9=5, 7", "35÷2=17, 1"),
    @("46÷5=9, 1", "95÷5=19, 0"),
    @("19÷6=3, 1", "94÷5=18, 4"),
    @("64÷9=7, 1", "89÷8=11, 1"),
    @("98÷6=16, 2", "92÷4=23, 0"),
    @("52÷7=7, 3", "91÷6=15, 1"),
    @("73÷8=9, 1", "31÷2=15, 1"),
    @("90÷6=15, 0", "57÷6=9, 3"),
    @("60÷4=15, 0", "10÷9=1, 1"),
    @("49÷7=7, 0", "54÷5=10, 4"),
    @("18÷9=2, 0", "85÷7=12, 1"),
    @("27÷5=5, 2", "54÷7=7, 5"),
    @("47÷9=5, 2", "56÷9=6, 2"),
    @("99÷2=49, 1", "55÷2=27, 1"),
    @("83÷6=13, 5", "62÷7=8, 6"),
    @("66÷4=16, 2", "74÷8=9, 2"),
    @("73÷4=18, 1", "51÷3=17, 0"),
    @("84÷2=42, 0", "28÷2=14, 0"),
    @("95÷4=23, 3", "80÷2=40, 0"),
    @("33÷7=4, 5", "69÷8=8, 5"),
    @("13÷9=1, 4", "71÷8=8, 7"),
    @("61÷5=12, 1", "54÷9=6, 0"),
    @("22÷6=3, 4", "88÷8=11, 0"),
    @("99÷6=16, 3", "16÷7=2, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
